# Update the "last source code" row (row 2) on Sheet1 with the new
# FuncLoc / SAID / Previous Doc values.
#
# Columns:
#   AV = FuncLoc       -> ABCD791616
#   AW = SAID          -> 9661848209
#   AX = Previous Doc  -> 5967521411
#
# AW2's cell is formatted with a numeric display format ("0") but must keep
# storing its value as text (matching how the sheet already stores similar
# all-digit identifiers as text elsewhere). Flipping the number format to
# Text ("@") before assigning the value -- and then restoring the original
# numeric format afterwards -- makes the engine keep the value as a string
# without altering the cell's final style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# FuncLoc (already alphanumeric, stays text automatically)
$ws.Range("AV2").Value = "ABCD791616"

# SAID - force text storage, then restore the original number format
$saidCell = $ws.Range("AW2")
$originalFormat = $saidCell.NumberFormat
$saidCell.NumberFormat = "@"
$saidCell.Value = "9661848209"
$saidCell.NumberFormat = $originalFormat

# Previous Doc (already alphanumeric-looking but all digits; cell format is
# text, so it is kept as text automatically)
$ws.Range("AX2").Value = "5967521411"
